$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the two unused sheets (Sheet2, Sheet3)
$wb.Worksheets.Item("Sheet2").Delete() | Out-Null
$wb.Worksheets.Item("Sheet3").Delete() | Out-Null

# Insert a new row at 60 (pushes the blank separator / totals rows down by one).
# Excel carries the row-59 cell styles down onto the fresh row automatically.
$ws.Range("A60:P60").Insert(-4121) | Out-Null

# Populate the new pub-run entry (White Hart, Duffield)
$ws.Range("A60").Value = 44097
$ws.Range("B60").Value = "The White Hart"
$ws.Range("C60").Value = "Duffield"
$ws.Range("D60").Value = "start/end at pub"
$ws.Range("E60").Value = 4.1399999999999997
$ws.Range("F60").Value = 0.043784722222222218
$ws.Range("G60").Formula = "=F60/E60"
$ws.Range("H60").Value = 1
$ws.Range("J60").Value = 1
$ws.Range("N60").Value = 1
$ws.Range("O60").Value = "Windley moated manorial complex"
$ws.Range("P60").Formula = "=SUM(H60:N60)*E60"

# Column J has no prior style to inherit (row 59 leaves it blank) - borrow it
# from the neighbouring H60 cell so it renders the same as the rest of the row
$ws.Range("H59").Copy() | Out-Null
$ws.Range("J60").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("J60").Value = 1

# Columns I, K, L and M stay blank in the new row - drop the leftover styled
# placeholder cells the insert carried over so the row matches columns
# A,B,C,D,E,F,G,H,J,N,O,P only
$ws.Range("I60").Clear() | Out-Null
$ws.Range("K60").Clear() | Out-Null
$ws.Range("L60").Clear() | Out-Null
$ws.Range("M60").Clear() | Out-Null

# Refresh the summary formulas on the totals row (now row 62) to include row 60
$ws.Range("E62").Formula = "=SUM(E5:E60)"
$ws.Range("G62").Formula = "=AVERAGE(G6:G60)"
$ws.Range("H62").Formula = "=SUM(H5:H60)"
$ws.Range("I62").Formula = "=SUM(I5:I59)"
$ws.Range("J62").Formula = "=SUM(J5:J59)"
$ws.Range("K62").Formula = "=SUM(K5:K59)"
$ws.Range("L62").Formula = "=SUM(L5:L59)"
$ws.Range("M62").Formula = "=SUM(M5:M59)"
$ws.Range("N62").Formula = "=SUM(N5:N60)"
$ws.Range("P62").Formula = "=SUM(P5:P61)"

# Restore the active selection to A60, as in the authored workbook
$ws.Range("A60").Select() | Out-Null
